# Apply changes described by the commit:
# "Versión 0.95 (15-09) cambios en los perfiles a mejorar"

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://core.hl7chile.cl/StructureDefinition/CodigoPaises"
$meta.Range("B4").Value = "PaisOrigen-Nacionalidad-Cl"
$meta.Range("B8").Value = "2021-09-15T12:30:30-03:00"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://core.hl7chile.cl/StructureDefinition/CodigoPaises"
$elements.Range("Y7").Value = "http://hl7.org/fhir/ValueSet/iso3166-1-N"

# Widen column Y to fit the new, longer value set URL (bestFit-style autofit).
# Target stored width is 37.8203125 characters; the ColumnWidth setter here
# snaps to whole-pixel increments (multiples of 1/6), so 37 is the closest
# input that lands on the nearest achievable grid point (37.8333...).
$elements.Columns.Item(25).ColumnWidth = 37
